# Update app version info (apps_info.xlsx) to the latest versions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Visual C Redistributable
$ws.Range("C2").Value = "0.33.0"

# Audacity
$ws.Range("C4").Value = "2.4.2"

# Calibre
$ws.Range("C5").Value = "4.20.0"

# cmder
$ws.Range("C6").Value = "1.3.15 2020/06/26"

# CrystalDiskInfo
$ws.Range("C7").Value = "8.7.0"

# Dropbox
$ws.Range("C9").Value = "101.4.434"

# Exiftool
$ws.Range("C13").Value = "12.01"

# Firefox
$ws.Range("C15").Value = "78.0.2"

# HWInfo
$ws.Range("C19").Value = "6.28"

# K-Lite Codec
$ws.Range("C21").Value = "15.6.0"

# MKVToolnix
$ws.Range("C22").Value = "48.0.0"

# Visual Studio Code
$ws.Range("C28").Value = "1.47"

# Row 30: Google Play Music Desktop -> YouTube Music (app replaced)
$ws.Range("A30").Value = "yt_music"
$ws.Range("B30").Value = "YouTube Music"
$ws.Range("C30").Value = "1.11.0"
$ws.Range("D30").Value = "https://github.com/ytmdesktop/ytmdesktop/releases"
$ws.Range("E30").Value = "https://github.com/ytmdesktop/ytmdesktop/releases"

# Java SE
$ws.Range("C31").Value = "more"

# Update the active cell selection
$null = $ws.Range("A2").Select()
